# 自动更新Excel文件 - decrement remaining-days column (E) by 1 for each
# data row, except row 36 (which is left untouched in the source diff),
# and special-case row 94 whose "剩余" (E) and "开始时间" (F) values were
# reset to new values rather than simply decremented.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column E = 5 ("剩余" / remaining), data rows run from 2 to 99.
for ($row = 2; $row -le 99; $row++) {
    if ($row -eq 36) {
        # Row 36 is unchanged in this update.
        continue
    }

    $cell = $ws.Cells.Item($row, 5)
    $current = $cell.Value2

    if ($row -eq 94) {
        # Row 94 gets an explicit reset rather than a simple decrement.
        $cell.Value2 = 7
        $ws.Cells.Item($row, 6).Value2 = 20251204
    } else {
        $cell.Value2 = $current - 1
    }
}
